$wb = $excel.ActiveWorkbook

# --- Sheet ALC, row 4 (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1072.381
$ws.Range("I4").Value = 1026.05
$ws.Range("J4").Value = 1999
$ws.Range("K4").Value = 1026.05
$ws.Range("L4").Value = 1999
$ws.Range("M4").Value = -912.05
$ws.Range("N4").Value = -2227

# --- Sheet ALC, row 33 (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 48915.047
$ws.Range("I33").Value = 62794.938
$ws.Range("J33").Value = 4499.4
$ws.Range("K33").Value = 62794.938
$ws.Range("L33").Value = 4499.4
$ws.Range("M33").Value = -62565.938
$ws.Range("N33").Value = -4957.4

# --- Sheet ALC, row 96 (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 621.375
$ws.Range("I96").Value = 621.375
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1864.125
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -491.125
$ws.Range("N96").ClearContents()

# --- Sheet ALC, row 98 (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1843.5714
$ws.Range("I98").Value = 2626.25
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 2626.25
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = -1128.25

# --- Sheet ALC, row 122 (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1843.5714
$ws.Range("I122").Value = 2626.25
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 7878.75
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -5428.75

# --- Sheet ALC, row 137 (hunk 5) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2019.1482
$ws.Range("I137").Value = 1676.9048
$ws.Range("J137").Value = 3217
$ws.Range("K137").Value = 5030.7144
$ws.Range("L137").Value = 9651
$ws.Range("M137").Value = -2480.7144
$ws.Range("N137").Value = -14751

# --- Sheet ARM, row 61 (hunk 6) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2556
$ws.Range("I61").Value = 1513.8572
$ws.Range("J61").Value = 3077.0715
$ws.Range("K61").Value = 1513.8572
$ws.Range("L61").Value = 3077.0715
$ws.Range("M61").Value = -1301.8572
$ws.Range("N61").Value = -3501.0715

# --- Sheet ARM, row 74 (hunk 7) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2263
$ws.Range("I74").Value = 1506.7826
$ws.Range("J74").Value = 4002.3
$ws.Range("K74").Value = 1506.7826
$ws.Range("L74").Value = 4002.3
$ws.Range("M74").Value = -632.7826
$ws.Range("N74").Value = -5750.3

# --- Sheet ARM, row 77 (hunk 8) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2263
$ws.Range("I77").Value = 1506.7826
$ws.Range("J77").Value = 4002.3
$ws.Range("K77").Value = 7533.913
$ws.Range("L77").Value = 20011.5
$ws.Range("M77").Value = -3165.913
$ws.Range("N77").Value = -28747.5

# --- Sheet ARM, row 136 (hunk 9) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2556
$ws.Range("I136").Value = 1513.8572
$ws.Range("J136").Value = 3077.0715
$ws.Range("K136").Value = 4541.571599999999
$ws.Range("L136").Value = 9231.2145
$ws.Range("M136").Value = -1991.571599999999
$ws.Range("N136").Value = -14331.2145

# --- Sheet BSM, row 20 (hunk 10) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 30836.143
$ws.Range("I20").Value = 42818.207
$ws.Range("J20").Value = 4693.4546
$ws.Range("K20").Value = 42818.207
$ws.Range("L20").Value = 4693.4546
$ws.Range("M20").Value = -42571.207
$ws.Range("N20").Value = -5187.4546

# --- Sheet BSM, row 134 (hunk 11) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2314.9778
$ws.Range("I134").Value = 2199
$ws.Range("J134").Value = 3242.8
$ws.Range("K134").Value = 6597
$ws.Range("L134").Value = 9728.400000000001
$ws.Range("M134").Value = -4062

# --- Sheet CRP, row 4 (hunk 12) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1490
$ws.Range("I4").Value = 1490
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1490
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1378

# --- Sheet CRP, row 7 (hunk 13) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 127.78571
$ws.Range("I7").Value = 45.22222
$ws.Range("J7").Value = 276.4
$ws.Range("K7").Value = 45.22222
$ws.Range("L7").Value = 276.4
$ws.Range("M7").Value = 67.77778000000001

# --- Sheet CRP, row 31 (hunk 14) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11703.768
$ws.Range("I31").Value = 13501.86
$ws.Range("J31").Value = 4601.3
$ws.Range("K31").Value = 13501.86
$ws.Range("L31").Value = 4601.3
$ws.Range("M31").Value = -13206.86
$ws.Range("N31").Value = -5191.3

# --- Sheet CRP, row 34 (hunk 15) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 11703.768
$ws.Range("I34").Value = 13501.86
$ws.Range("J34").Value = 4601.3
$ws.Range("K34").Value = 13501.86
$ws.Range("L34").Value = 4601.3
$ws.Range("M34").Value = -13299.86
$ws.Range("N34").Value = -5005.3

# --- Sheet CRP, row 58 (hunk 16) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9608
$ws.Range("I58").Value = 1311.8334
$ws.Range("J58").Value = 38052
$ws.Range("K58").Value = 1311.8334
$ws.Range("L58").Value = 38052
$ws.Range("M58").Value = -1108.8334
$ws.Range("N58").Value = -38458

# --- Sheet CRP, row 132 (hunk 17) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2548.121
$ws.Range("I132").Value = 1871.4615
$ws.Range("J132").Value = 5061.4287
$ws.Range("K132").Value = 5614.3845
$ws.Range("L132").Value = 15184.2861
$ws.Range("M132").Value = -3084.3845
$ws.Range("N132").Value = -20244.2861

# --- Sheet CRP, row 134 (hunk 18) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1091.0625
$ws.Range("I134").Value = 1091.0625
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3273.1875
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -738.1875
$ws.Range("N134").ClearContents()

# --- Sheet CRP, row 136 (hunk 19) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 9608
$ws.Range("I136").Value = 1311.8334
$ws.Range("J136").Value = 38052
$ws.Range("K136").Value = 3935.5002
$ws.Range("L136").Value = 114156
$ws.Range("M136").Value = -1385.5002
$ws.Range("N136").Value = -119256

# --- Sheet CUL, row 120 (hunk 20) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 4053.3333
$ws.Range("I120").Value = 4053.3333
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 12159.9999
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -7321.999899999999
$ws.Range("N120").ClearContents()

# --- Sheet GSM, row 116 (hunk 21) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 58000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 58000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 58000
$ws.Range("N116").Value = -67178

# --- Sheet GSM, row 119 (hunk 22) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 40000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 40000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676

# --- Sheet GSM, row 124 (hunk 23) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 43750
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 43750
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 43750
$ws.Range("N124").Value = -53570

# --- Sheet GSM, row 126 (hunk 24) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2548.8262
$ws.Range("I126").Value = 3040.4
$ws.Range("J126").Value = 2170.6924
$ws.Range("K126").Value = 9121.200000000001
$ws.Range("L126").Value = 6512.0772
$ws.Range("M126").Value = -6651.200000000001
$ws.Range("N126").Value = -11452.0772

# --- Sheet LTW, row 132 (hunk 25) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5609.143
$ws.Range("I132").Value = 5664.353
$ws.Range("J132").Value = 5374.5
$ws.Range("K132").Value = 16993.059
$ws.Range("L132").Value = 16123.5
$ws.Range("M132").Value = -14463.059
$ws.Range("N132").Value = -21183.5

# --- Sheet LTW, row 136 (hunk 26) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2127.52
$ws.Range("I136").Value = 2085.3333
$ws.Range("J136").Value = 2349
$ws.Range("K136").Value = 6255.999899999999
$ws.Range("L136").Value = 7047
$ws.Range("M136").Value = -3705.999899999999
$ws.Range("N136").Value = -12147

# --- Sheet WVR, row 132 (hunk 27) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2680.05
$ws.Range("I132").Value = 2342
$ws.Range("J132").Value = 3187.125
$ws.Range("K132").Value = 7026
$ws.Range("L132").Value = 9561.375
$ws.Range("M132").Value = -4496
$ws.Range("N132").Value = -14621.375

# --- Sheet WVR, row 136 (hunk 28) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1307.561
$ws.Range("I136").Value = 468.16666
$ws.Range("J136").Value = 1964.4783
$ws.Range("K136").Value = 1404.49998
$ws.Range("L136").Value = 5893.4349
$ws.Range("M136").Value = 1145.50002
$ws.Range("N136").Value = -10993.4349

